$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the "Periodo Mora" rows (16-19) into ascending period order,
# carrying each period's "Valor Mora" (column F) along with it.
$ws.Range("E16").Value = "1805"
$ws.Range("E17").Value = "1806"
$ws.Range("E18").Value = "1807"
$ws.Range("E19").Value = "1903"

$ws.Range("F16").Value = 31249
$ws.Range("F17").Value = 31249
$ws.Range("F18").Value = 31249
$ws.Range("F19").Value = 33125
